$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRoutes = @(
    @('AURANGABAD','Unknown','WayPoint','N19°51''39.95"','E075°24''18.90"','Unknown Name'),
    @('DUBOX','Unknown','WayPoint','N20°22''15.00"','E076°49''56.00"','Unknown Name'),
    @('NINAT','Unknown','WayPoint','N20°32''48.00"','E077°24''56.00"','Unknown Name'),
    @('NAGPUR','Unknown','WayPoint','N21°04''53.22"','E079°03''22.67"','Unknown Name'),
    @('TEGIG','Unknown','WayPoint','N21°30''09.00"','E080°35''37.00"','Unknown Name'),
    @('DOSAT','Unknown','WayPoint','N21°41''32.99"','E081°18''24.00"','Unknown Name'),
    @('OPONI','Unknown','WayPoint','N21°47''03.00"','E081°38''24.00"','Unknown Name'),
    @('NIPAD','Unknown','WayPoint','N21°52''17.99"','E081°59''52.99"','Unknown Name'),
    @('OTABA','Unknown','WayPoint','N22°11''02.99"','E083°19''22.99"','Unknown Name'),
    @('KINKI','Unknown','WayPoint','N22°19''17.99"','E083°55''38.00"','Unknown Name'),
    @('AGROM','Unknown','WayPoint','N22°31''44.99"','E084°49''59.99"','Unknown Name'),
    @('JAMSHEDPUR','Unknown','WayPoint','N22°48''47.47"','E086°10''26.15"','Unknown Name'),
    @('KAMGU','Unknown','WayPoint','N14°58''31.00"','E079°25''12.00"','Unknown Name'),
    @('KIKUR','Unknown','WayPoint','N15°58''54.99"','E079°10''32.00"','Unknown Name'),
    @('ANDAV','Unknown','WayPoint','N17°28''18.00"','E078°48''49.00"','Unknown Name'),
    @('ALBED','Unknown','WayPoint','N18°07''01.99"','E078°39''23.99"','Unknown Name'),
    @('LAROB','Unknown','WayPoint','N19°10''14.00"','E078°23''40.00"','Unknown Name'),
    @('SUDEL','Unknown','WayPoint','N19°34''17.99"','E078°17''18.99"','Unknown Name'),
    @('TASEX','Unknown','WayPoint','N20°33''54.00"','E078°01''36.99"','Unknown Name'),
    @('RENAG','Unknown','WayPoint','N22°03''12.00"','E077°39''32.00"','Unknown Name'),
    @('BHOPAL','Unknown','WayPoint','N23°16''59.73"','E077°20''11.85"','Unknown Name'),
    @('INTIL','Unknown','WayPoint','N26°27''51.00"','E076°32''40.99"','Unknown Name'),
)

$startRow = 346
for ($i = 0; $i -lt $newRoutes.Count; $i++) {
    $r = $startRow + $i
    $row = $newRoutes[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
